$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 38), columns E:H ---
# Reuse the bordered header style already used by row 2 (C2:G2) via copy/paste
# of formatting only, so the existing style index (s="2") is reused instead of
# minting a new one.
$ws.Range("C2").Copy()
$ws.Range("E38:H38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E38").Value = "システム内客数（計測値）"
$ws.Range("F38").Value = "システム内時間（計測値）"
$ws.Range("G38").Value = "システム内客数（理論値）"
$ws.Range("H38").Value = "システム内時間（理論値）"

# --- Row 39: システムA ---
$ws.Range("D39").Value = "システムA"
$ws.Range("E39").Value = 0.17668600000000001
$ws.Range("F39").Value = 0.59756900000000002
$ws.Range("G39").Value = 0.2
$ws.Range("H39").Value = 0.6

# --- Row 40: システムB ---
$ws.Range("D40").Value = "システムB"
$ws.Range("E40").Value = 0.507081
$ws.Range("F40").Value = 0.50390400000000002

# --- Row 41: システムC ---
$ws.Range("D41").Value = "システムC"
$ws.Range("E41").Value = 0.19773499999999999
$ws.Range("F41").Value = 0.200348
$ws.Range("G41").Value = 0.2
$ws.Range("H41").Value = 0.2

# Match the final selection recorded in the saved workbook.
$ws.Range("K41").Select()
